$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.542.70'
$ws.Range("E2").Value = '  +0.70%  '
$ws.Range("D3").Value = '3.250.61'
$ws.Range("E3").Value = '  +2.78%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '607.57'
$ws.Range("E5").Value = '  +0.94%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '157.34'
$ws.Range("E6").Value = '  +2.19%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '3.251.62'
$ws.Range("E8").Value = '  +2.87%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.550'
$ws.Range("E9").Value = '  +0.41%  '
$ws.Range("E10").Value = '  +2.67%  '
$ws.Range("E11").Value = '  +5.89%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.501'
$ws.Range("E12").Value = '  -3.42%  '
$ws.Range("E13").Value = '  +1.82%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '39.11'
$ws.Range("E14").Value = '  +2.18%  '
$ws.Range("D15").Value = '3.789.45'
$ws.Range("E15").Value = '  +2.95%  '
$ws.Range("D16").Value = '66.679.81'
$ws.Range("E16").Value = '  +0.87%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.46'
$ws.Range("E17").Value = '  +0.79%  '
$ws.Range("D18").Value = '3.256.80'
$ws.Range("E18").Value = '  +2.91%  '
$ws.Range("E19").Value = '  +1.14%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '506.50'
$ws.Range("E20").Value = '  -0.61%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.43'
$ws.Range("E21").Value = '  +0.43%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.752'
$ws.Range("E22").Value = '  +3.46%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.13'
$ws.Range("E23").Value = '  +1.37%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '14.76'
$ws.Range("E24").Value = '  +0.99%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '87.09'
$ws.Range("E25").Value = '  +3.13%  '
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.03'
$ws.Range("E27").Value = '  +1.56%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.16'
$ws.Range("E28").Value = '  +1.40%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.41'
$ws.Range("E29").Value = '  +1.17%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.133'
$ws.Range("E30").Value = '  +50.56%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.91'
$ws.Range("E31").Value = '  -4.54%  '
$ws.Range("E32").Value = '  -0.68%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '28.03'
$ws.Range("E33").Value = '  +0.36%  '
$ws.Range("E34").Value = '  +0.21%  '
$ws.Range("E35").Value = '  -2.96%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.47'
$ws.Range("E36").Value = '  -0.10%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.37'
$ws.Range("E37").Value = '  +20.61%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '55.63'
$ws.Range("E38").Value = '  +1.85%  '
$ws.Range("D39").Value = '0.0₃0782'
$ws.Range("E39").Value = '  +16.38%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '494.78'
$ws.Range("E40").Value = '  +0.07%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0425'
$ws.Range("E41").Value = '  +1.43%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.128'
$ws.Range("E42").Value = '  -0.51%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.83'
$ws.Range("E43").Value = '  +1.01%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.294'
$ws.Range("E44").Value = '  -0.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.52'
$ws.Range("E45").Value = '  +4.71%  '
$ws.Range("D46").Value = '2.987.68'
$ws.Range("E46").Value = '  +6.07%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '28.93'
$ws.Range("E47").Value = '  +4.18%  '
$ws.Range("E48").Value = '  +5.99%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.119'
$ws.Range("E49").Value = '  +2.50%  '
$ws.Range("E50").Value = '  -0.05%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '121.02'
$ws.Range("E51").Value = '  -0.43%  '
